# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) for the rows whose figures moved in this update.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the workbook's existing text-formatted numbers/
# percentages) instead of auto-converting them into numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.93"
$ws.Range("E2").Value = "'0.20%"
$ws.Range("D3").Value = "'41.59"
$ws.Range("E3").Value = "'0.98%"
$ws.Range("D4").Value = "'5.691"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("D5").Value = "'0.08386"
$ws.Range("E5").Value = "'3.92%"
$ws.Range("D6").Value = "'8.813"
$ws.Range("E6").Value = "'0.78%"
$ws.Range("D7").Value = "'1.999"
$ws.Range("E7").Value = "'-1.56%"
$ws.Range("D8").Value = "'4.479"
$ws.Range("E8").Value = "'-0.99%"
$ws.Range("D10").Value = "'0.9239"
$ws.Range("E10").Value = "'0.29%"
$ws.Range("D11").Value = "'0.1276"
$ws.Range("E11").Value = "'1.47%"
$ws.Range("D12").Value = "'0.1963"
$ws.Range("E12").Value = "'1.00%"
$ws.Range("D13").Value = "'0.09364"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("D14").Value = "'0.03892"
$ws.Range("E14").Value = "'4.98%"
$ws.Range("E15").Value = "'0.57%"
$ws.Range("D16").Value = "'0.001303"
$ws.Range("E16").Value = "'0.45%"
$ws.Range("D17").Value = "'0.006111"
$ws.Range("E17").Value = "'-2.09%"
$ws.Range("E18").Value = "'1.86%"
$ws.Range("D19").Value = "'0.3500"
$ws.Range("E19").Value = "'0.69%"
$ws.Range("D20").Value = "'8.936"
$ws.Range("E20").Value = "'7.51%"
$ws.Range("D21").Value = "'0.1364"
$ws.Range("E21").Value = "'-3.74%"
$ws.Range("D22").Value = "'0.2511"
$ws.Range("E22").Value = "'-5.47%"
$ws.Range("D23").Value = "'0.04417"
$ws.Range("E23").Value = "'-0.39%"
$ws.Range("D24").Value = "'0.001244"
$ws.Range("E24").Value = "'-1.26%"
$ws.Range("D25").Value = "'0.004385"
$ws.Range("E25").Value = "'1.86%"
$ws.Range("D26").Value = "'0.0001191"
$ws.Range("E26").Value = "'-4.22%"
$ws.Range("D27").Value = "'0.0003994"
$ws.Range("E27").Value = "'0.03%"
$ws.Range("D39").Value = "'0.02817"
$ws.Range("E39").Value = "'-1.62%"
$ws.Range("D40").Value = "'0.05524"
$ws.Range("E40").Value = "'1.14%"
$ws.Range("D41").Value = "'0.007948"
$ws.Range("E41").Value = "'2.03%"
$ws.Range("D42").Value = "'0.1436"
$ws.Range("E42").Value = "'1.32%"
$ws.Range("D43").Value = "'0.008975"
$ws.Range("E43").Value = "'-10.04%"
$ws.Range("D44").Value = "'0.002092"
$ws.Range("E44").Value = "'-6.32%"
$ws.Range("D45").Value = "'0.01176"
$ws.Range("E45").Value = "'-0.99%"
$ws.Range("D46").Value = "'0.00006946"
$ws.Range("E46").Value = "'2.43%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.18%"
$ws.Range("D48").Value = "'0.003304"
$ws.Range("E48").Value = "'9.39%"
$ws.Range("D49").Value = "'0.002280"
$ws.Range("E49").Value = "'-0.20%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.18%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.18%"
